# Apply updated crypto price/volume data to Sheet1 (columns D and E).
# Values that look like plain decimal numbers (e.g. "1.005") are written
# with NumberFormat "@" first so Excel keeps them as text, matching the
# original inline-string cell contents (prices like "25.888.99" or
# percentages like "  -0.87%  " are unambiguous text and need no such hint).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '25.888.99'
# Row 3: Ethereum
$ws.Range("D3").Value = '1.637.58'
$ws.Range("E3").Value = '  -0.87%  '
# Row 4: TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.05%  '
# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.22'
$ws.Range("E5").Value = '  -0.03%  '
# Row 6: XRP
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5023'
$ws.Range("E6").Value = '  -1.81%  '
# Row 7: USDC
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.002'
$ws.Range("E7").Value = '  -0.19%  '
# Row 8: Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2565'
$ws.Range("E8").Value = '  -0.75%  '
# Row 9: Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06387'
$ws.Range("E9").Value = '  -0.72%  '
# Row 10: Solana
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.70'
$ws.Range("E10").Value = '  -1.22%  '
# Row 11: TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07726'
$ws.Range("E11").Value = '  -0.82%  '
# Row 12: WrappedEther
$ws.Range("D12").Value = '1.654.34'
$ws.Range("E12").Value = '  +0.18%  '
# Row 13: Polkadot
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.269'
$ws.Range("E13").Value = '  -0.28%  '
# Row 14: WrappedliquidstakedEther2.0
$ws.Range("D14").Value = '1.862.03'
$ws.Range("E14").Value = '  -0.90%  '
# Row 15: Polygon
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5448'
$ws.Range("E15").Value = '  -1.20%  '
# Row 16: ShibaInu
$ws.Range("D16").Value = '0.0₅7907'
$ws.Range("E16").Value = '  -1.21%  '
# Row 17: Litecoin
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.26'
$ws.Range("E17").Value = '  +0.23%  '
# Row 18: WrappedBTC
$ws.Range("D18").Value = '25.884.09'
$ws.Range("E18").Value = '  -1.06%  '
# Row 19: Dai
$ws.Range("E19").Value = '  -0.12%  '
# Row 20: BitcoinCash
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '202.85'
$ws.Range("E20").Value = '  -3.74%  '
# Row 21: Uniswap
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.379'
$ws.Range("E21").Value = '  -0.30%  '
# Row 22: Avalanche
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.913'
$ws.Range("E22").Value = '  -1.52%  '
# Row 23: Chainlink
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.973'
$ws.Range("E23").Value = '  -1.31%  '
# Row 24: BinanceUSD
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.003'
$ws.Range("E24").Value = '  -0.12%  '
# Row 25: Toncoin
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.931'
$ws.Range("E25").Value = '  +10.48%  '
# Row 26: Monero
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '141.46'
$ws.Range("E26").Value = '  -1.76%  '
# Row 27: Stellar
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1137'
$ws.Range("E27").Value = '  -3.38%  '
# Row 28: EthereumClassic
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.68'
$ws.Range("E28").Value = '  -0.67%  '
# Row 29: Cosmos
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.728'
$ws.Range("E29").Value = '  -3.59%  '
# Row 30: PancakeSwap
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.242'
$ws.Range("E30").Value = '  +0.17%  '
# Row 31: Hedera
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.04931'
$ws.Range("E31").Value = '  -4.09%  '
# Row 32: InternetComputer(DFINITY)
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.273'
$ws.Range("E32").Value = '  -2.13%  '
# Row 33: Filecoin
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.193'
$ws.Range("E33").Value = '  -0.69%  '
# Row 34: LidoDAOToken
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.542'
$ws.Range("E34").Value = '  -1.01%  '
# Row 35: HuobiToken
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.368'
$ws.Range("E35").Value = '  +0.52%  '
# Row 36: MXToken
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.633'
$ws.Range("E36").Value = '  -3.81%  '
# Row 37: ARBITRUM
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.8926'
$ws.Range("E37").Value = '  -3.42%  '
# Row 38: Maker
$ws.Range("D38").Value = '1.156.90'
$ws.Range("E38").Value = '  -0.93%  '
# Row 39: ImmutableX
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5596'
$ws.Range("E39").Value = '  -1.73%  '
# Row 40: VeChain
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01562'
$ws.Range("E40").Value = '  -1.37%  '
# Row 41: PaxDollar
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.001'
$ws.Range("E41").Value = '  -0.30%  '
# Row 42: FraxShare
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.706'
$ws.Range("E42").Value = '  +0.86%  '
# Row 43: TrustWalletToken
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8085'
$ws.Range("E43").Value = '  -1.97%  '
# Row 44: Quant
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.56'
$ws.Range("E44").Value = '  -0.73%  '
# Row 45: RocketPoolETH
$ws.Range("D45").Value = '1.774.59'
$ws.Range("E45").Value = '  -0.82%  '
# Row 46: BabyDogeCoin
$ws.Range("E46").Value = '  -0.07%  '
# Row 47: Mantle
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4527'
$ws.Range("E47").Value = '  -0.47%  '
# Row 48: Frax
$ws.Range("E48").Value = '  -0.25%  '
# Row 49: Aave
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '55.02'
$ws.Range("E49").Value = '  -0.77%  '
# Row 50: Cronos
$ws.Range("E50").Value = '  -0.03%  '
# Row 51: USDD
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.003'
$ws.Range("E51").Value = '  -0.34%  '
